$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CustomerMappingDriver Class section ---
# Row 29: "For successfully scanning data from input file"
#   score drops from 16 to 15, and a grading comment is added.
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = "For not adding scanned products to inventory"

# Row 30: "For correct and properly aligned output"
#   grading comment reworded to reflect compilation errors causing no output.
$ws.Range("F30").Value = "(-4) For no output for all methods due to compilation errors"

# --- Generic section ---
# Row 37: "Compilation errors if any"
#   grading comment reworded with more detail about which classes failed to compile.
$ws.Range("F37").Value = "(-5) For compilation errors in driver as well as CustomerMapping class"

# Leave the final selection on the last edited cell, matching the author's
# last position in the worksheet when they saved.
$ws.Range("F37").Select()
